$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:U61")
try {
    $lo = $ws.ListObjects.Add(1, $rng, $null, 1, $null, "")
    Write-Host "added with empty style name"
} catch {
    Write-Host "ERR add: $_"
}
